$d = $word.ActiveDocument

# Locate the "Architecture Flow" heading paragraph
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Architecture Flow*") {
        $target = $p
    }
}

# Insert a brand new empty paragraph right after it, styled as Source Code
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Style = "SourceCode"

$lines = @(
    'graph TD',
    '    A[Local Raw Data <br> transactions.csv] -->|Mounted Volume| B(Airflow Worker Container)',
    '    ',
    '    subgraph Docker runtime via Colima',
    '        B -->|Pandas Transformation| C{process_data.py}',
    '        C -->|Outputs Processed Data| D[Local Processed Data <br> processed_transactions.csv]',
    '    end',
    '    ',
    '    D -->|upload_to_gcs.py| E[(GCS Data Lake)]',
    '    E -->|load_to_bq.py <br> Append Mode| F[(BigQuery Data Warehouse)]',
    '    ',
    '    G((Airflow Scheduler)) -.->|Triggers DAG Daily| B',
    '    ',
    '    classDef storage fill:#f9f,stroke:#333,stroke-width:2px;',
    '    classDef script fill:#bbf,stroke:#333,stroke-width:2px;',
    '    classDef cloud fill:#ff9,stroke:#333,stroke-width:2px;',
    '    ',
    '    class A,D storage;',
    '    class C script;',
    '    class E,F cloud;'
)

for ($i = 0; $i -lt $lines.Length; $i++) {
    $line = $lines[$i]

    # The paragraph's range always ends right after its (still unwritten-to)
    # paragraph mark, so "End - 1" is exactly where the next run should land.
    $lineStart = $newPara.Range.End - 1
    $newPara.Range.InsertAfter($line)

    $lineRange = $d.Range($lineStart, $lineStart + $line.Length)
    $lineRange.Style = "VerbatimChar"

    if ($i -lt ($lines.Length - 1)) {
        $breakPos = $lineStart + $line.Length
        $breakPoint = $d.Range($breakPos, $breakPos)
        $breakPoint.InsertBreak(6)
    }
}

Write-Output "done"
